$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 13.172

$ws.Range("C4").Value = -12.835
$ws.Range("D4").Value = -8.168000000000001
$ws.Range("E4").Value = 12.647

$ws.Range("D5").Value = -8.581999999999999

$ws.Range("C6").Value = -12.445

$ws.Range("C7").Value = -13.276

$ws.Range("D8").Value = -8.216999999999999

$ws.Range("E9").Value = 12.961

$ws.Range("E11").Value = 13.035

$ws.Range("E14").Value = 13.06

$ws.Range("C16").Value = -12.072
$ws.Range("D16").Value = -8.625999999999999

$ws.Range("E18").Value = 12.84

$ws.Range("C20").Value = -13.041

$ws.Range("D22").Value = -8.147

$ws.Range("E25").Value = 13.201
